$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four holiday rows currently all say "TAMIL PONGAL" in column A.
# Replace them with the four distinct Pongal-day names (Bhogi, Thai,
# Mattu, Kannum) used by the fixed Tamil Nadu holiday template.
$ws.Range("A2").Value = "BHOGI PONGAL"
$ws.Range("A3").Value = "THAI PONGAL"
$ws.Range("A4").Value = "MATTU PONGAL"
$ws.Range("A5").Value = "KANNUM PONGAL"

# Move the active selection back to the top of the (now empty) data
# entry area instead of leaving it down at C21:C22.
$ws.Range("A6").Select()
